# Sheet1: add two "anchor" numeric cells (B1, A2) styled bold/centered with
# a thin box border, and a label cell (B2) holding the shared string
# "disconnected_elements".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the styled look on B1 first ...
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4160     # xlTop
$ws.Range("B1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("B1").Borders.Weight = 2            # xlThin

# ... then clone the resulting style onto A2 via copy/paste-format so both
# cells end up sharing the very same cell-format record.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = $false
